$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-11-15 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-16 Thursday", 2) | Out-Null

# Update each cell of the 20x5 answer table with its new value
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "49+19=68"
$t.Cell(1,2).Range.Text = "20-4=16"
$t.Cell(1,3).Range.Text = "71-36=35"
$t.Cell(1,4).Range.Text = "42-3=39"
$t.Cell(1,5).Range.Text = "6+66=72"
$t.Cell(2,1).Range.Text = "8+74=82"
$t.Cell(2,2).Range.Text = "94-47=47"
$t.Cell(2,3).Range.Text = "6+17=23"
$t.Cell(2,4).Range.Text = "15+57=72"
$t.Cell(2,5).Range.Text = "68+24=92"
$t.Cell(3,1).Range.Text = "80-27=53"
$t.Cell(3,2).Range.Text = "16+46=62"
$t.Cell(3,3).Range.Text = "40-2=38"
$t.Cell(3,4).Range.Text = "70-27=43"
$t.Cell(3,5).Range.Text = "62-7=55"
$t.Cell(4,1).Range.Text = "23+29=52"
$t.Cell(4,2).Range.Text = "56-29=27"
$t.Cell(4,3).Range.Text = "94-48=46"
$t.Cell(4,4).Range.Text = "69+2=71"
$t.Cell(4,5).Range.Text = "40-24=16"
$t.Cell(5,1).Range.Text = "69+3=72"
$t.Cell(5,2).Range.Text = "47+8=55"
$t.Cell(5,3).Range.Text = "79+18=97"
$t.Cell(5,4).Range.Text = "27+26=53"
$t.Cell(5,5).Range.Text = "95-39=56"
$t.Cell(6,1).Range.Text = "81-36=45"
$t.Cell(6,2).Range.Text = "5+46=51"
$t.Cell(6,3).Range.Text = "51-18=33"
$t.Cell(6,4).Range.Text = "81-32=49"
$t.Cell(6,5).Range.Text = "94-26=68"
$t.Cell(7,1).Range.Text = "80-29=51"
$t.Cell(7,2).Range.Text = "59+23=82"
$t.Cell(7,3).Range.Text = "73-15=58"
$t.Cell(7,4).Range.Text = "77+7=84"
$t.Cell(7,5).Range.Text = "35-19=16"
$t.Cell(8,1).Range.Text = "94-49=45"
$t.Cell(8,2).Range.Text = "6+8=14"
$t.Cell(8,3).Range.Text = "27+17=44"
$t.Cell(8,4).Range.Text = "74-16=58"
$t.Cell(8,5).Range.Text = "16-7=9"
$t.Cell(9,1).Range.Text = "44+27=71"
$t.Cell(9,2).Range.Text = "17+76=93"
$t.Cell(9,3).Range.Text = "53-27=26"
$t.Cell(9,4).Range.Text = "80-73=7"
$t.Cell(9,5).Range.Text = "14+28=42"
$t.Cell(10,1).Range.Text = "69+22=91"
$t.Cell(10,2).Range.Text = "54-25=29"
$t.Cell(10,3).Range.Text = "53-45=8"
$t.Cell(10,4).Range.Text = "25+18=43"
$t.Cell(10,5).Range.Text = "71-8=63"
$t.Cell(11,1).Range.Text = "9+12=21"
$t.Cell(11,2).Range.Text = "66-17=49"
$t.Cell(11,3).Range.Text = "20-14=6"
$t.Cell(11,4).Range.Text = "51-27=24"
$t.Cell(11,5).Range.Text = "92-28=64"
$t.Cell(12,1).Range.Text = "77+19=96"
$t.Cell(12,2).Range.Text = "91-9=82"
$t.Cell(12,3).Range.Text = "46+19=65"
$t.Cell(12,4).Range.Text = "52-44=8"
$t.Cell(12,5).Range.Text = "91-53=38"
$t.Cell(13,1).Range.Text = "41-37=4"
$t.Cell(13,2).Range.Text = "25+36=61"
$t.Cell(13,3).Range.Text = "30-18=12"
$t.Cell(13,4).Range.Text = "60-14=46"
$t.Cell(13,5).Range.Text = "54-6=48"
$t.Cell(14,1).Range.Text = "32-27=5"
$t.Cell(14,2).Range.Text = "49+22=71"
$t.Cell(14,3).Range.Text = "79+15=94"
$t.Cell(14,4).Range.Text = "78-69=9"
$t.Cell(14,5).Range.Text = "97-69=28"
$t.Cell(15,1).Range.Text = "77-58=19"
$t.Cell(15,2).Range.Text = "52-6=46"
$t.Cell(15,3).Range.Text = "67+27=94"
$t.Cell(15,4).Range.Text = "39+45=84"
$t.Cell(15,5).Range.Text = "28+48=76"
$t.Cell(16,1).Range.Text = "94-89=5"
$t.Cell(16,2).Range.Text = "27+6=33"
$t.Cell(16,3).Range.Text = "48+45=93"
$t.Cell(16,4).Range.Text = "48-9=39"
$t.Cell(16,5).Range.Text = "40-35=5"
$t.Cell(17,1).Range.Text = "29+42=71"
$t.Cell(17,2).Range.Text = "63-16=47"
$t.Cell(17,3).Range.Text = "48-19=29"
$t.Cell(17,4).Range.Text = "59+24=83"
$t.Cell(17,5).Range.Text = "97-49=48"
$t.Cell(18,1).Range.Text = "81-65=16"
$t.Cell(18,2).Range.Text = "18+77=95"
$t.Cell(18,3).Range.Text = "78+18=96"
$t.Cell(18,4).Range.Text = "87+6=93"
$t.Cell(18,5).Range.Text = "79+3=82"
$t.Cell(19,1).Range.Text = "85-76=9"
$t.Cell(19,2).Range.Text = "24+37=61"
$t.Cell(19,3).Range.Text = "28+33=61"
$t.Cell(19,4).Range.Text = "18+4=22"
$t.Cell(19,5).Range.Text = "38+9=47"
$t.Cell(20,1).Range.Text = "24+39=63"
$t.Cell(20,2).Range.Text = "31-26=5"
$t.Cell(20,3).Range.Text = "8+7=15"
$t.Cell(20,4).Range.Text = "26+27=53"
$t.Cell(20,5).Range.Text = "62-33=29"
